$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, pushing the current rows 6-16 down to 7-17
$ws.Rows("6:6").Insert()

# Populate the new row 6 with a new weekly price record (mirrors the other
# rows' static fields, with its own date/volume/price values)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C6").Value = "Los Lagos"
$ws.Range("D6").Value = 44789
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 100112035
$ws.Range("G6").Value = "Bruselas (repollito)"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 90
$ws.Range("K6").Value = 24000
$ws.Range("L6").Value = 24000
$ws.Range("M6").Value = 24000
$ws.Range("N6").Value = "$/malla 15 kilos"
$ws.Range("O6").Value = "Provincia de Quillota"
$ws.Range("P6").Value = 1600
$ws.Range("Q6").Value = 15
$ws.Range("R6").Value = "Hortaliza"
